# Update the prediction-score column ("1-g__Barnesiella") with the
# freshly computed values from the latest ful-path.csv run, replacing
# the previous placeholder values of 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 5879.8783442027498
$ws.Range("B3").Value = 3757.1413770939603
$ws.Range("B4").Value = 2490.3190413447705
